# Change the table style on every table in the deck that currently uses the
# old custom "Table_0" style GUID, switching it to the built-in style GUID
# {A62A2CE7-E7AF-4133-A4F9-6FD590943077}.
#
# (Three tables - on the slides holding gridCol widths 2879725 / 3424250 /
# 2881325 EMU, i.e. slides 14, 15 and 16 - use the old style
# {1267FB70-1E0E-49CD-B6B1-96DCA16BB781}.)

$OldStyleId = "{1267FB70-1E0E-49CD-B6B1-96DCA16BB781}"
$NewStyleId = "{A62A2CE7-E7AF-4133-A4F9-6FD590943077}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $tbl = $shape.Table

            if ($tbl.Style -eq $OldStyleId) {
                $tbl.ApplyStyle($NewStyleId)
            }
        }
    }
}
